# Add data for 2022-04-05
# Updates:
#  - Sheet name "Through 2022-03-27" -> "Through 2022-03-28"
#  - Header label "2022 (through 03-27)" -> "2022 (through 03-28)"
#  - I4 (April Total) 111 -> 120
#  - I14 (Grand Total) 411 -> 420

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-03-28"

$ws.Range("I1").Value = "2022 (through 03-28)"

$ws.Range("I4").Value = 120
$ws.Range("I14").Value = 420
